$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 237 (old rows 237-241 shift down to 239-243)
$ws.Rows.Item(237).Insert()
$ws.Rows.Item(238).Insert()

# Copy cell formatting (styles) from the template row 236 into the two new rows
$ws.Range("A236:AC236").Copy()
$ws.Range("A237:AC238").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 237
$ws.Cells.Item(237,1).Value = 235
$ws.Cells.Item(237,2).Value = 6775582
$ws.Cells.Item(237,3).Value = "Poland Ekstraklasa"
$ws.Cells.Item(237,4).Value = "Poland Ekstraklasa"
$ws.Cells.Item(237,5).Value = 45388.41666666666
$ws.Cells.Item(237,6).Value = "Korona Kielce"
$ws.Cells.Item(237,7).Value = "Stal Mielec"
$ws.Cells.Item(237,8).Value = 1
$ws.Cells.Item(237,9).Value = 0
$ws.Cells.Item(237,10).Value = "H"
$ws.Cells.Item(237,11).Value = 1.8
$ws.Cells.Item(237,12).Value = 3.5
$ws.Cells.Item(237,13).Value = 4.5
$ws.Cells.Item(237,14).Value = 1.85
$ws.Cells.Item(237,15).Value = 3.4
$ws.Cells.Item(237,16).Value = 4.333
$ws.Cells.Item(237,17).Value = -0.5
$ws.Cells.Item(237,18).Value = 1.875
$ws.Cells.Item(237,19).Value = 1.975
$ws.Cells.Item(237,20).Value = 2.25
$ws.Cells.Item(237,21).Value = 1.95
$ws.Cells.Item(237,22).Value = 1.9
$ws.Cells.Item(237,23).Value = 0.8500000000000001
$ws.Cells.Item(237,24).Value = -1
$ws.Cells.Item(237,25).Value = -1
$ws.Cells.Item(237,26).Value = 0.875
$ws.Cells.Item(237,27).Value = -1
$ws.Cells.Item(237,28).Value = -1
$ws.Cells.Item(237,29).Value = 0.8999999999999999

# Row 238
$ws.Cells.Item(238,1).Value = 236
$ws.Cells.Item(238,2).Value = 6774876
$ws.Cells.Item(238,3).Value = "Poland Ekstraklasa"
$ws.Cells.Item(238,4).Value = "Poland Ekstraklasa"
$ws.Cells.Item(238,5).Value = 45388.52083333334
$ws.Cells.Item(238,6).Value = "Ruch Chorzow"
$ws.Cells.Item(238,7).Value = "Puszcza Niepolomice"
$ws.Cells.Item(238,8).Value = 0
$ws.Cells.Item(238,9).Value = 0
$ws.Cells.Item(238,10).Value = "D"
$ws.Cells.Item(238,11).Value = 1.85
$ws.Cells.Item(238,12).Value = 3.5
$ws.Cells.Item(238,13).Value = 4.2
$ws.Cells.Item(238,14).Value = 1.8
$ws.Cells.Item(238,15).Value = 3.6
$ws.Cells.Item(238,16).Value = 4.2
$ws.Cells.Item(238,17).Value = -0.5
$ws.Cells.Item(238,18).Value = 1.85
$ws.Cells.Item(238,19).Value = 2
$ws.Cells.Item(238,20).Value = 2.5
$ws.Cells.Item(238,21).Value = 2
$ws.Cells.Item(238,22).Value = 1.85
$ws.Cells.Item(238,23).Value = -1
$ws.Cells.Item(238,24).Value = 2.6
$ws.Cells.Item(238,25).Value = -1
$ws.Cells.Item(238,26).Value = -1
$ws.Cells.Item(238,27).Value = 1
$ws.Cells.Item(238,28).Value = -1
$ws.Cells.Item(238,29).Value = 0.8500000000000001

# Row 239
$ws.Cells.Item(239,1).Value = 237
$ws.Cells.Item(239,2).Value = 6775587
$ws.Cells.Item(239,3).Value = "Poland Ekstraklasa"
$ws.Cells.Item(239,4).Value = "Poland Ekstraklasa"
$ws.Cells.Item(239,5).Value = 45388.625
$ws.Cells.Item(239,6).Value = "Slask Wroclaw"
$ws.Cells.Item(239,7).Value = "Warta Poznan"
$ws.Cells.Item(239,8).Value = 2
$ws.Cells.Item(239,9).Value = 1
$ws.Cells.Item(239,10).Value = "H"
$ws.Cells.Item(239,11).Value = 1.727
$ws.Cells.Item(239,12).Value = 3.8
$ws.Cells.Item(239,13).Value = 4.5
$ws.Cells.Item(239,14).Value = 1.909
$ws.Cells.Item(239,15).Value = 3
$ws.Cells.Item(239,16).Value = 5
$ws.Cells.Item(239,17).Value = -0.5
$ws.Cells.Item(239,18).Value = 1.975
$ws.Cells.Item(239,19).Value = 1.875
$ws.Cells.Item(239,20).Value = 2
$ws.Cells.Item(239,21).Value = 2.05
$ws.Cells.Item(239,22).Value = 1.8
$ws.Cells.Item(239,23).Value = 0.909
$ws.Cells.Item(239,24).Value = -1
$ws.Cells.Item(239,25).Value = -1
$ws.Cells.Item(239,26).Value = 0.9750000000000001
$ws.Cells.Item(239,27).Value = -1
$ws.Cells.Item(239,28).Value = 1.05
$ws.Cells.Item(239,29).Value = -1

# Row 240
$ws.Cells.Item(240,1).Value = 238
$ws.Cells.Item(240,2).Value = 6775586
$ws.Cells.Item(240,3).Value = "Poland Ekstraklasa"
$ws.Cells.Item(240,4).Value = "Poland Ekstraklasa"
$ws.Cells.Item(240,5).Value = 45389.3125
$ws.Cells.Item(240,6).Value = "Widzew Lodz"
$ws.Cells.Item(240,7).Value = "Piast Gliwice"
$ws.Cells.Item(240,8).Value = 1
$ws.Cells.Item(240,9).Value = 0
$ws.Cells.Item(240,10).Value = "H"
$ws.Cells.Item(240,11).Value = 2.5
$ws.Cells.Item(240,12).Value = 3.25
$ws.Cells.Item(240,13).Value = 2.8
$ws.Cells.Item(240,14).Value = 2.25
$ws.Cells.Item(240,15).Value = 3.2
$ws.Cells.Item(240,16).Value = 3.25
$ws.Cells.Item(240,17).Value = -0.25
$ws.Cells.Item(240,18).Value = 1.95
$ws.Cells.Item(240,19).Value = 1.9
$ws.Cells.Item(240,20).Value = 2
$ws.Cells.Item(240,21).Value = 1.9
$ws.Cells.Item(240,22).Value = 1.95
$ws.Cells.Item(240,23).Value = 1.25
$ws.Cells.Item(240,24).Value = -1
$ws.Cells.Item(240,25).Value = -1
$ws.Cells.Item(240,26).Value = 0.95
$ws.Cells.Item(240,27).Value = -1
$ws.Cells.Item(240,28).Value = -1
$ws.Cells.Item(240,29).Value = 0.95

# Row 241
$ws.Cells.Item(241,1).Value = 239
$ws.Cells.Item(241,2).Value = 6775583
$ws.Cells.Item(241,3).Value = "Poland Ekstraklasa"
$ws.Cells.Item(241,4).Value = "Poland Ekstraklasa"
$ws.Cells.Item(241,5).Value = 45389.41666666666
$ws.Cells.Item(241,6).Value = "Lech Poznan"
$ws.Cells.Item(241,7).Value = "Pogon Szczecin"
$ws.Cells.Item(241,8).Value = 1
$ws.Cells.Item(241,9).Value = 0
$ws.Cells.Item(241,10).Value = "H"
$ws.Cells.Item(241,11).Value = 2.5
$ws.Cells.Item(241,12).Value = 3.4
$ws.Cells.Item(241,13).Value = 2.7
$ws.Cells.Item(241,14).Value = 2.2
$ws.Cells.Item(241,15).Value = 3.5
$ws.Cells.Item(241,16).Value = 3.1
$ws.Cells.Item(241,17).Value = -0.25
$ws.Cells.Item(241,18).Value = 1.95
$ws.Cells.Item(241,19).Value = 1.9
$ws.Cells.Item(241,20).Value = 2.75
$ws.Cells.Item(241,21).Value = 2
$ws.Cells.Item(241,22).Value = 1.85
$ws.Cells.Item(241,23).Value = 1.2
$ws.Cells.Item(241,24).Value = -1
$ws.Cells.Item(241,25).Value = -1
$ws.Cells.Item(241,26).Value = 0.95
$ws.Cells.Item(241,27).Value = -1
$ws.Cells.Item(241,28).Value = -1
$ws.Cells.Item(241,29).Value = 0.8500000000000001

# Row 242
$ws.Cells.Item(242,1).Value = 240
$ws.Cells.Item(242,2).Value = 6775584
$ws.Cells.Item(242,3).Value = "Poland Ekstraklasa"
$ws.Cells.Item(242,4).Value = "Poland Ekstraklasa"
$ws.Cells.Item(242,5).Value = 45389.52083333334
$ws.Cells.Item(242,6).Value = "Legia Warsaw"
$ws.Cells.Item(242,7).Value = "Jagiellonia Bialystok"
$ws.Cells.Item(242,8).Value = 1
$ws.Cells.Item(242,9).Value = 1
$ws.Cells.Item(242,10).Value = "D"
$ws.Cells.Item(242,11).Value = 2
$ws.Cells.Item(242,12).Value = 3.5
$ws.Cells.Item(242,13).Value = 3.6
$ws.Cells.Item(242,14).Value = 1.909
$ws.Cells.Item(242,15).Value = 3.5
$ws.Cells.Item(242,16).Value = 3.8
$ws.Cells.Item(242,17).Value = -0.5
$ws.Cells.Item(242,18).Value = 2
$ws.Cells.Item(242,19).Value = 1.85
$ws.Cells.Item(242,20).Value = 2.75
$ws.Cells.Item(242,21).Value = 1.95
$ws.Cells.Item(242,22).Value = 1.9
$ws.Cells.Item(242,23).Value = -1
$ws.Cells.Item(242,24).Value = 2.5
$ws.Cells.Item(242,25).Value = -1
$ws.Cells.Item(242,26).Value = -1
$ws.Cells.Item(242,27).Value = 0.8500000000000001
$ws.Cells.Item(242,28).Value = -1
$ws.Cells.Item(242,29).Value = 0.8999999999999999

# Row 243
$ws.Cells.Item(243,1).Value = 241
$ws.Cells.Item(243,2).Value = 6775579
$ws.Cells.Item(243,3).Value = "Poland Ekstraklasa"
$ws.Cells.Item(243,4).Value = "Poland Ekstraklasa"
$ws.Cells.Item(243,5).Value = 45390.58333333334
$ws.Cells.Item(243,6).Value = "Zaglebie Lubin"
$ws.Cells.Item(243,7).Value = "Gornik Zabrze"
$ws.Cells.Item(243,11).Value = 2.2
$ws.Cells.Item(243,12).Value = 3.4
$ws.Cells.Item(243,13).Value = 3.2
$ws.Cells.Item(243,14).Value = 2.3
$ws.Cells.Item(243,15).Value = 3.3
$ws.Cells.Item(243,16).Value = 3.1
$ws.Cells.Item(243,17).Value = -0.25
$ws.Cells.Item(243,18).Value = 2.025
$ws.Cells.Item(243,19).Value = 1.825
$ws.Cells.Item(243,20).Value = 2.5
$ws.Cells.Item(243,21).Value = 1.975
$ws.Cells.Item(243,22).Value = 1.875
$ws.Cells.Item(243,23).Value = 0
$ws.Cells.Item(243,24).Value = 0
$ws.Cells.Item(243,25).Value = 0
$ws.Cells.Item(243,26).Value = 0
$ws.Cells.Item(243,27).Value = 0

